$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.636.64"
$ws.Range("E2").Value = "  -3.07%  "
$ws.Range("D3").Value = "2.610.81"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'573.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.05%  "
$ws.Range("D6").Value = "'154.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.629"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  -5.51%  "
$ws.Range("D10").Value = "'5.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").Value = "'0.384"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.06%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "'28.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").Value = "3.081.04"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").Value = "'0.0000182"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.18%  "
$ws.Range("D16").Value = "63.531.91"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "2.620.50"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "'12.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.42%  "
$ws.Range("E19").Value = "  -2.42%  "
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").Value = "'342.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'66.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.92%  "
$ws.Range("D24").Value = "'1.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("D25").Value = "'0.0000107"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.49%  "
$ws.Range("D26").Value = "'9.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.08%  "
$ws.Range("D27").Value = "'576.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.09%  "
$ws.Range("D28").Value = "'1.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("D31").Value = "'7.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").Value = "'2.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.27%  "
$ws.Range("E33").Value = "  -3.62%  "
$ws.Range("D34").Value = "'6.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").Value = "'5.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").Value = "'0.408"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("D37").Value = "'19.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'153.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "'41.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("D43").Value = "'156.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("D44").Value = "'2.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.65%  "
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("D46").Value = "'22.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "'0.0592"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.627"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").Value = "'0.0250"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").Value = "'18.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.34%  "
